$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing date-format style (bold, bordered, centered, custom
# date/time number format) from the last pre-existing data row (A385) down
# across the newly appended rows in column A, same as the source workbook
# does for every row in column A.
$ws.Range("A385").Copy() | Out-Null
$ws.Range("A386:A464").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$data = @(
  @(44460,0,0,0),
  @(44461,0,0,0),
  @(44462,0,0,0),
  @(44463,0,0,0),
  @(44464,0,0,0),
  @(44465,0,0,0),
  @(44466,0,0,0),
  @(44467,0,0,0),
  @(44468,0,0,0),
  @(44469,0,0,0),
  @(44470,0,0,0),
  @(44471,0,0,0),
  @(44472,0,0,0),
  @(44473,0,0,0),
  @(44474,0,0,0),
  @(44475,0,0,0),
  @(44476,0,0,0),
  @(44477,0,0,0),
  @(44478,0,0,0),
  @(44479,0,0,0),
  @(44480,0,0,0),
  @(44481,0,0,0),
  @(44482,0,0,0),
  @(44483,0,0,0),
  @(44484,0,0,0),
  @(44485,0,0,0),
  @(44486,0,0,0),
  @(44487,0,0,0),
  @(44488,0,0,0),
  @(44489,0,0,0),
  @(44490,0,0,0),
  @(44491,0,0,0),
  @(44492,0,0,0),
  @(44493,0,0,0),
  @(44494,0,0,0),
  @(44495,0,0,0),
  @(44496,0,0,0),
  @(44497,2,2,96.15384615384616),
  @(44498,0,2,96.15384615384616),
  @(44499,0,2,96.15384615384616),
  @(44500,0,2,96.15384615384616),
  @(44501,0,2,96.15384615384616),
  @(44502,1,3,144.2307692307692),
  @(44503,0,3,144.2307692307692),
  @(44504,0,1,48.07692307692308),
  @(44505,0,1,48.07692307692308),
  @(44506,0,1,48.07692307692308),
  @(44507,0,1,48.07692307692308),
  @(44508,0,1,48.07692307692308),
  @(44509,0,0,0),
  @(44510,0,0,0),
  @(44511,0,0,0),
  @(44512,0,0,0),
  @(44513,0,0,0),
  @(44514,0,0,0),
  @(44515,0,0,0),
  @(44516,1,1,48.07692307692308),
  @(44517,0,1,48.07692307692308),
  @(44518,0,1,48.07692307692308),
  @(44519,1,2,96.15384615384616),
  @(44520,0,2,96.15384615384616),
  @(44521,0,2,96.15384615384616),
  @(44522,0,2,96.15384615384616),
  @(44523,0,1,48.07692307692308),
  @(44524,0,1,48.07692307692308),
  @(44525,1,2,96.15384615384616),
  @(44526,1,2,96.15384615384616),
  @(44527,1,3,144.2307692307692),
  @(44528,0,3,144.2307692307692),
  @(44529,0,3,144.2307692307692),
  @(44530,0,3,144.2307692307692),
  @(44531,0,3,144.2307692307692),
  @(44532,0,2,96.15384615384616),
  @(44533,0,1,48.07692307692308),
  @(44534,0,0,0),
  @(44535,0,0,0),
  @(44536,0,0,0),
  @(44537,1,1,48.07692307692308),
  @(44538,0,1,48.07692307692308)
)
$startRow = 386
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
    $ws.Cells.Item($r, 4).Value = $vals[3]
}
